$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Market Cap) values for existing rows 2-9
$ws.Range("C2").Value = 709599331094.9697
$ws.Range("C3").Value = 232009741643.4652
$ws.Range("C4").Value = 36763362405.97121
$ws.Range("C5").Value = 32303492845.14954
$ws.Range("C6").Value = 23863430866.82423
$ws.Range("C7").Value = 12826871086.8618
$ws.Range("C8").Value = 11677096606.99001
$ws.Range("C9").Value = 9088763660.586021

# Row 10 and 11 swap: Polygon (MATIC-USD) and Toncoin (TON-USD) swap places
$ws.Range("A10").Value = "Toncoin"
$ws.Range("B10").Value = "TON-USD"
$ws.Range("C10").Value = 7814775280.991043

$ws.Range("A11").Value = "Polygon"
$ws.Range("B11").Value = "MATIC-USD"
$ws.Range("C11").Value = 7539009532.292098

# Rows 12-22 values update
$ws.Range("C12").Value = 7533781949.363321
$ws.Range("C13").Value = 7408357493.749443
$ws.Range("C14").Value = 6602397511.290756
$ws.Range("C15").Value = 5958607945.875978
$ws.Range("C16").Value = 5111727373.900389
$ws.Range("C17").Value = 5062042475.280413
$ws.Range("C18").Value = 4431124989.914836
$ws.Range("C19").Value = 3770368267.006455
$ws.Range("C20").Value = 3456637581.265154
$ws.Range("C21").Value = 3350340934.019464
$ws.Range("C22").Value = 3259399994.280393

# Rows 23-25 rotate: Kaspa -> Monero -> Ethereum Classic -> Kaspa
$ws.Range("A23").Value = "Monero"
$ws.Range("B23").Value = "XMR-USD"
$ws.Range("C23").Value = 2984873211.79517

$ws.Range("A24").Value = "Ethereum Classic"
$ws.Range("B24").Value = "ETC-USD"
$ws.Range("C24").Value = 2726804834.093751

$ws.Range("A25").Value = "Kaspa"
$ws.Range("B25").Value = "KAS-USD"
$ws.Range("C25").Value = 2634899949.649961

# Row 26
$ws.Range("C26").Value = 2427815139.415774
